$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("measles_costs")

# Update vacc_cost input (B14)
$ws.Range("B14").Value = 50

# Updated contact/case numbers in column E (rows 19-38)
$ws.Range("E19").Value = 17583
$ws.Range("E20").Value = 4417
$ws.Range("E21").Value = 13325
$ws.Range("E22").Value = 10242
$ws.Range("E23").Value = 18421
$ws.Range("E24").Value = 3619
$ws.Range("E25").Value = 4267
$ws.Range("E26").Value = 2798
$ws.Range("E27").Value = 4496
$ws.Range("E28").Value = 2257
$ws.Range("E29").Value = 2942
$ws.Range("E30").Value = 854
$ws.Range("E31").Value = 8153
$ws.Range("E32").Value = 1317
$ws.Range("E33").Value = 2804
$ws.Range("E34").Value = 11016
$ws.Range("E35").Value = 689
$ws.Range("E36").Value = 16845
$ws.Range("E37").Value = 660
$ws.Range("E38").Value = 1328

# Update sheet view: scroll position and selection
$ws.Range("A18:W38").Select()
$excel.ActiveWindow.ScrollRow = 6
$excel.ActiveWindow.ScrollColumn = 4
